$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5; existing rows 5-47 shift down to 6-48
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with this week's data (same constant columns as other rows)
$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(5, 3).Value = "Coquimbo"
$ws.Cells.Item(5, 4).Value = 44699
$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 100112026
$ws.Cells.Item(5, 7).Value = "Haba"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 400
$ws.Cells.Item(5, 11).Value = 14000
$ws.Cells.Item(5, 12).Value = 15000
$ws.Cells.Item(5, 13).Value = 14500
$ws.Cells.Item(5, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(5, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(5, 16).Value = 580
$ws.Cells.Item(5, 17).Value = 25
$ws.Cells.Item(5, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
